$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.428.72"
$ws.Range("E2").Value = "  +1.85%  "

# Row 3
$ws.Range("D3").Value = "3.290.44"
$ws.Range("E3").Value = "  +1.56%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.88%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.47"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.53%  "

# Row 7
$ws.Range("E7").Value = "  -0.07%  "

# Row 8
$ws.Range("E8").Value = "  +0.58%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.136"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.67%  "

# Row 10
$ws.Range("E10").Value = "  -0.07%  "

# Row 11
$ws.Range("E11").Value = "  +2.08%  "

# Row 12
$ws.Range("D12").Value = "3.860.75"
$ws.Range("E12").Value = "  +1.52%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.137"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.03%  "

# Row 15
$ws.Range("D15").Value = "68.413.32"
$ws.Range("E15").Value = "  +1.92%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000173"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.15%  "

# Row 17
$ws.Range("D17").Value = "3.257.68"
$ws.Range("E17").Value = "  +0.27%  "

# Row 18
$ws.Range("E18").Value = "  +1.27%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "384.52"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.11%  "

# Row 21
$ws.Range("E21").Value = "  +3.18%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.52"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.52%  "

# Row 23
$ws.Range("E23").Value = "  -0.11%  "

# Row 24
$ws.Range("B24").Value = "PEPE"
$ws.Range("C24").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000123"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.78%  "

# Row 25
$ws.Range("B25").Value = "Polygon"
$ws.Range("C25").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.517"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.73%  "

# Row 26
$ws.Range("E26").Value = "  +8.13%  "

# Row 27
$ws.Range("E27").Value = "  +2.01%  "

# Row 28
$ws.Range("E28").Value = "  -0.56%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.78%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.59%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.06"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.16%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.32"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.90%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.24"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.51%  "

# Row 34
$ws.Range("E34").Value = "  +0.03%  "

# Row 35
$ws.Range("E35").Value = "  +3.48%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.29%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.89"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.61%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.838"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.46%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.87"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.81%  "

# Row 40
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.68"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.00%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.75"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.32%  "

# Row 42
$ws.Range("E42").Value = "  +5.32%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.52%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0693"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.14%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.33"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.46%  "

# Row 46
$ws.Range("D46").Value = "2.638.88"
$ws.Range("E46").Value = "  -3.08%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "341.03"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -6.33%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0285"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.68%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "32.31"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.51%  "

# Row 50
$ws.Range("E50").Value = "  +1.77%  "

# Row 51
$ws.Range("E51").Value = "  -0.12%  "
